# Update sample TSV: add "Refrigerator (4 Celsius)" to the
# specimen_preservation_temperature list (#575).
#
# This inserts a new allowed value "Refrigerator (4 Celsius)" into the
# "specimen_prese...mperature list" lookup sheet (between the existing
# "Freezer (-20 Celsius)" and "Room Temperature" entries), extends the
# dimension/validation range on that sheet from A1:A5 to A1:A6, and updates
# the data validation on the main "Export as TSV" sheet's column L
# (specimen_preservation_temperature) to reference the larger range and use
# the generic "Value must come from ... list." error message.

$wb = $excel.ActiveWorkbook

# --- 1. Update the lookup list sheet ------------------------------------
$listWs = $wb.Worksheets.Item("specimen_prese...mperature list")

# Preserve the existing last entry ("Room Temperature") by pushing it down
# to row 6, then insert the new value in the now-vacated row 5 -- this
# matches the add-new-option-before-the-last-one ordering used upstream.
$lastValue = $listWs.Range("A5").Value()
$listWs.Range("A6").Value = $lastValue
$listWs.Range("A5").Value = "Refrigerator (4 Celsius)"

# --- 2. Update the data validation on the main sheet --------------------
$mainWs = $wb.Worksheets.Item("Export as TSV")
$colRange = $mainWs.Range("L2:L1048576")
$validation = $colRange.Validation

$xlValidateList = 3
$xlValidAlertStop = 1
$xlBetween = 1

$validation.Modify($xlValidateList, $xlValidAlertStop, $xlBetween, "'specimen_prese...mperature list'!`$A`$1:`$A`$6")
$validation.ErrorMessage = "Value must come from specimen_prese...mperature list."
